$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H19").Value = 937.2857
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 937.2857
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 937.2857
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -1287.2857
$ws.Range("H38").Value = 1058.75
$ws.Range("I38").Value = 956.8421
$ws.Range("K38").Value = 2870.5263
$ws.Range("M38").Value = -2498.5263
$ws.Range("H40").Value = 2845.6667
$ws.Range("I40").Value = 2881.6365
$ws.Range("K40").Value = 2881.6365
$ws.Range("M40").Value = -2706.6365
$ws.Range("H41").Value = 528.4545000000001
$ws.Range("I41").Value = 826.6667
$ws.Range("J41").Value = 322
$ws.Range("K41").Value = 826.6667
$ws.Range("L41").Value = 322
$ws.Range("M41").Value = -386.6667
$ws.Range("N41").Value = -1202
$ws.Range("H48").Value = 900
$ws.Range("I48").Value = 801
$ws.Range("J48").Value = 999
$ws.Range("K48").Value = 2403
$ws.Range("L48").Value = 2997
$ws.Range("M48").Value = -2111
$ws.Range("N48").Value = -3581
$ws.Range("H56").Value = 900
$ws.Range("I56").Value = 801
$ws.Range("J56").Value = 999
$ws.Range("K56").Value = 2403
$ws.Range("L56").Value = 2997
$ws.Range("M56").Value = -1869
$ws.Range("N56").Value = -4065
$ws.Range("H64").Value = 4852.2
$ws.Range("I64").Value = 4500.3335
$ws.Range("K64").Value = 4500.3335
$ws.Range("M64").Value = -4252.3335
$ws.Range("H67").Value = 4852.2
$ws.Range("I67").Value = 4500.3335
$ws.Range("K67").Value = 4500.3335
$ws.Range("M67").Value = -3642.3335
$ws.Range("H70").Value = 1273.75
$ws.Range("I70").Value = 1265.6666
$ws.Range("K70").Value = 3796.9998
$ws.Range("M70").Value = -3526.9998
$ws.Range("H73").Value = 1273.75
$ws.Range("I73").Value = 1265.6666
$ws.Range("K73").Value = 3796.9998
$ws.Range("M73").Value = -2860.9998
$ws.Range("H76").Value = 5099.8
$ws.Range("I76").Value = 4833
$ws.Range("J76").Value = 5500
$ws.Range("K76").Value = 4833
$ws.Range("L76").Value = 5500
$ws.Range("M76").Value = -4518
$ws.Range("N76").Value = -6130
$ws.Range("H79").Value = 5099.8
$ws.Range("I79").Value = 4833
$ws.Range("J79").Value = 5500
$ws.Range("K79").Value = 4833
$ws.Range("L79").Value = 5500
$ws.Range("M79").Value = -3741
$ws.Range("N79").Value = -7684
$ws.Range("H86").Value = 4615.1816
$ws.Range("I86").Value = 1618.5555
$ws.Range("K86").Value = 1618.5555
$ws.Range("M86").Value = -495.5554999999999
$ws.Range("H89").Value = 4615.1816
$ws.Range("I89").Value = 1618.5555
$ws.Range("K89").Value = 8092.7775
$ws.Range("M89").Value = -2476.7775
$ws.Range("H99").Value = 1714.875
$ws.Range("J99").Value = 4154.6665
$ws.Range("L99").Value = 12463.9995
$ws.Range("N99").Value = -15459.9995
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H106").Value = 1488.9286
$ws.Range("I106").Value = 1518.9
$ws.Range("K106").Value = 1518.9
$ws.Range("M106").Value = -887.9000000000001
$ws.Range("H127").Value = 2420.25
$ws.Range("I127").Value = 2420.25
$ws.Range("K127").Value = 7260.75
$ws.Range("M127").Value = -2300.75
$ws.Range("H131").Value = 4111.6665
$ws.Range("I131").Value = 1667.5
$ws.Range("J131").Value = 9000
$ws.Range("K131").Value = 5002.5
$ws.Range("L131").Value = 27000
$ws.Range("M131").Value = 37.5
$ws.Range("N131").Value = -37080
$ws.Range("H132").Value = 112781.555
$ws.Range("I132").Value = 112781.555
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 338344.665
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -335814.665
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 595.7826
$ws.Range("I135").Value = 395.5909
$ws.Range("K135").Value = 3560.3181
$ws.Range("M135").Value = -1025.3181
$ws.Range("H137").Value = 1563.1724
$ws.Range("I137").Value = 1510
$ws.Range("K137").Value = 4530
$ws.Range("M137").Value = -1980
$ws.Range("H138").Value = 2481.1143
$ws.Range("J138").Value = 2827
$ws.Range("L138").Value = 8481
$ws.Range("N138").Value = -18761
$ws.Range("H141").Value = 3629.2415
$ws.Range("I141").Value = 3702.1155
$ws.Range("K141").Value = 11106.3465
$ws.Range("M141").Value = -5926.3465

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1335.2941
$ws.Range("I2").Value = 1447.7858
$ws.Range("J2").Value = 810.3333
$ws.Range("K2").Value = 1447.7858
$ws.Range("L2").Value = 810.3333
$ws.Range("M2").Value = -1334.7858
$ws.Range("N2").Value = -1036.3333
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 1300.2703
$ws.Range("I32").Value = 1300.2703
$ws.Range("K32").Value = 1300.2703
$ws.Range("M32").Value = -1013.2703
$ws.Range("H61").Value = 3871.1667
$ws.Range("I61").Value = 3305.0908
$ws.Range("K61").Value = 3305.0908
$ws.Range("M61").Value = -3093.0908
$ws.Range("H102").Value = 2409.5
$ws.Range("I102").Value = 2409.5
$ws.Range("K102").Value = 2409.5
$ws.Range("M102").Value = -787.5
$ws.Range("H110").Value = 2870.6667
$ws.Range("I110").Value = 1806
$ws.Range("K110").Value = 1806
$ws.Range("M110").Value = 239
$ws.Range("H116").Value = 1335.2941
$ws.Range("I116").Value = 1447.7858
$ws.Range("J116").Value = 810.3333
$ws.Range("K116").Value = 1447.7858
$ws.Range("L116").Value = 810.3333
$ws.Range("M116").Value = 846.2141999999999
$ws.Range("N116").Value = -5398.3333
$ws.Range("H132").Value = 21743212
$ws.Range("I132").Value = 2150.6843
$ws.Range("K132").Value = 6452.0529
$ws.Range("M132").Value = -3922.0529
$ws.Range("H136").Value = 3871.1667
$ws.Range("I136").Value = 3305.0908
$ws.Range("K136").Value = 9915.2724
$ws.Range("M136").Value = -7365.2724

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1335.2941
$ws.Range("I3").Value = 1447.7858
$ws.Range("J3").Value = 810.3333
$ws.Range("K3").Value = 1447.7858
$ws.Range("L3").Value = 810.3333
$ws.Range("M3").Value = -1333.7858
$ws.Range("N3").Value = -1038.3333
$ws.Range("H20").Value = 849.3
$ws.Range("I20").Value = 858.6111
$ws.Range("K20").Value = 858.6111
$ws.Range("M20").Value = -611.6111
$ws.Range("H86").Value = 4449.9707
$ws.Range("I86").Value = 2283.611
$ws.Range("K86").Value = 2283.611
$ws.Range("M86").Value = -1160.611
$ws.Range("H89").Value = 4449.9707
$ws.Range("I89").Value = 2283.611
$ws.Range("K89").Value = 11418.055
$ws.Range("M89").Value = -5802.055
$ws.Range("H99").Value = 782.8333
$ws.Range("I99").Value = 782.8333
$ws.Range("K99").Value = 782.8333
$ws.Range("M99").Value = 715.1667
$ws.Range("H105").Value = 2487.5
$ws.Range("I105").Value = 2098.6086
$ws.Range("K105").Value = 2098.6086
$ws.Range("M105").Value = -351.6086
$ws.Range("H134").Value = 38892316
$ws.Range("I134").Value = 22730056
$ws.Range("K134").Value = 68190168
$ws.Range("M134").Value = -68187633

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1287.5385
$ws.Range("I16").Value = 985.2727
$ws.Range("K16").Value = 985.2727
$ws.Range("M16").Value = -698.2727
$ws.Range("H22").Value = 358
$ws.Range("I22").Value = 350
$ws.Range("J22").Value = 360.66666
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 360.66666
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = -1060.66666
$ws.Range("H31").Value = 2340
$ws.Range("I31").Value = 2234.0715
$ws.Range("K31").Value = 2234.0715
$ws.Range("M31").Value = -1939.0715
$ws.Range("H34").Value = 2340
$ws.Range("I34").Value = 2234.0715
$ws.Range("K34").Value = 2234.0715
$ws.Range("M34").Value = -2032.0715
$ws.Range("H58").Value = 1528.2439
$ws.Range("I58").Value = 1627.9286
$ws.Range("K58").Value = 1627.9286
$ws.Range("M58").Value = -1424.9286
$ws.Range("H99").Value = 1864.3636
$ws.Range("I99").Value = 1694.1428
$ws.Range("J99").Value = 2162.25
$ws.Range("K99").Value = 1694.1428
$ws.Range("L99").Value = 2162.25
$ws.Range("M99").Value = -196.1428000000001
$ws.Range("N99").Value = -5158.25
$ws.Range("H113").Value = 1287.5385
$ws.Range("I113").Value = 985.2727
$ws.Range("K113").Value = 985.2727
$ws.Range("M113").Value = 1184.7273
$ws.Range("H126").Value = 1864.3636
$ws.Range("I126").Value = 1694.1428
$ws.Range("J126").Value = 2162.25
$ws.Range("K126").Value = 5082.428400000001
$ws.Range("L126").Value = 6486.75
$ws.Range("M126").Value = -2612.428400000001
$ws.Range("N126").Value = -11426.75
$ws.Range("H132").Value = 10651.28
$ws.Range("I132").Value = 11287.85
$ws.Range("K132").Value = 33863.55
$ws.Range("M132").Value = -31333.55
$ws.Range("H134").Value = 5884977.5
$ws.Range("I134").Value = 2472.5715
$ws.Range("K134").Value = 7417.7145
$ws.Range("M134").Value = -4882.7145
$ws.Range("H136").Value = 1528.2439
$ws.Range("I136").Value = 1627.9286
$ws.Range("K136").Value = 4883.7858
$ws.Range("M136").Value = -2333.7858

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 550
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value = 550
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H56").Value = 7614.278
$ws.Range("I56").Value = 7614.278
$ws.Range("K56").Value = 7614.278
$ws.Range("M56").Value = -7084.278
$ws.Range("H80").Value = 4135.2856
$ws.Range("I80").Value = 3486.75
$ws.Range("K80").Value = 10460.25
$ws.Range("M80").Value = -9524.25
$ws.Range("H83").Value = 4135.2856
$ws.Range("I83").Value = 3486.75
$ws.Range("K83").Value = 31380.75
$ws.Range("M83").Value = -26700.75
$ws.Range("H88").Value = 2999
$ws.Range("I88").Value = 2999
$ws.Range("K88").Value = 8997
$ws.Range("M88").Value = -8569
$ws.Range("H91").Value = 2999
$ws.Range("I91").Value = 2999
$ws.Range("K91").Value = 8997
$ws.Range("M91").Value = -7515
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H100").Value = 16509.334
$ws.Range("J100").Value = 16509.334
$ws.Range("L100").Value = 49528.00199999999
$ws.Range("N100").Value = -51150.00199999999
$ws.Range("H107").Value = 832.4706
$ws.Range("J107").Value = 857.6
$ws.Range("L107").Value = 2572.8
$ws.Range("N107").Value = -6412.8
$ws.Range("H116").Value = 119421.164
$ws.Range("I116").Value = 140105.4
$ws.Range("K116").Value = 420316.2
$ws.Range("M116").Value = -416874.2
$ws.Range("H132").Value = 849
$ws.Range("J132").Value = 698.5
$ws.Range("L132").Value = 6286.5
$ws.Range("N132").Value = -11346.5
$ws.Range("H137").Value = 398504.7
$ws.Range("I137").Value = 3853.6
$ws.Range("J137").Value = 539451.5
$ws.Range("K137").Value = 11560.8
$ws.Range("L137").Value = 1618354.5
$ws.Range("M137").Value = -6460.799999999999
$ws.Range("N137").Value = -1628554.5

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7116.7036
$ws.Range("I70").Value = 7039.579
$ws.Range("J70").Value = 7299.875
$ws.Range("K70").Value = 7039.579
$ws.Range("L70").Value = 7299.875
$ws.Range("M70").Value = -6769.579
$ws.Range("N70").Value = -7839.875
$ws.Range("H73").Value = 7116.7036
$ws.Range("I73").Value = 7039.579
$ws.Range("J73").Value = 7299.875
$ws.Range("K73").Value = 7039.579
$ws.Range("L73").Value = 7299.875
$ws.Range("M73").Value = -6103.579
$ws.Range("N73").Value = -9171.875
$ws.Range("H80").Value = 5114.5
$ws.Range("I80").Value = 3181.25
$ws.Range("J80").Value = 6081.125
$ws.Range("K80").Value = 3181.25
$ws.Range("L80").Value = 6081.125
$ws.Range("M80").Value = -2183.25
$ws.Range("N80").Value = -8077.125
$ws.Range("H83").Value = 5114.5
$ws.Range("I83").Value = 3181.25
$ws.Range("J83").Value = 6081.125
$ws.Range("K83").Value = 15906.25
$ws.Range("L83").Value = 30405.625
$ws.Range("M83").Value = -10914.25
$ws.Range("N83").Value = -40389.625
$ws.Range("H102").Value = 1302.1578
$ws.Range("I102").Value = 978.75
$ws.Range("J102").Value = 1537.3636
$ws.Range("K102").Value = 978.75
$ws.Range("L102").Value = 1537.3636
$ws.Range("M102").Value = 643.25
$ws.Range("N102").Value = -4781.3636
$ws.Range("H126").Value = 6465.643
$ws.Range("I126").Value = 14002.667
$ws.Range("J126").Value = 4410.091
$ws.Range("K126").Value = 42008.001
$ws.Range("L126").Value = 13230.273
$ws.Range("M126").Value = -39538.001
$ws.Range("N126").Value = -18170.273
$ws.Range("H132").Value = 1988.6666
$ws.Range("I132").Value = 2026.2
$ws.Range("J132").Value = 1801
$ws.Range("K132").Value = 6078.6
$ws.Range("L132").Value = 5403
$ws.Range("M132").Value = -3548.6
$ws.Range("N132").Value = -10463

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2998.3333
$ws.Range("I7").Value = 2998.5
$ws.Range("K7").Value = 2998.5
$ws.Range("M7").Value = -2886.5
$ws.Range("H46").Value = 3672.3157
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3672.3157
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3672.3157
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4048.3157
$ws.Range("H61").Value = 1969.4546
$ws.Range("J61").Value = 2781.5
$ws.Range("L61").Value = 2781.5
$ws.Range("N61").Value = -3185.5
$ws.Range("H82").Value = 1573.125
$ws.Range("I82").Value = 1626.4286
$ws.Range("K82").Value = 1626.4286
$ws.Range("M82").Value = -1265.4286
$ws.Range("H85").Value = 1573.125
$ws.Range("I85").Value = 1626.4286
$ws.Range("K85").Value = 1626.4286
$ws.Range("M85").Value = -378.4286
$ws.Range("H109").Value = 44998.5
$ws.Range("J109").Value = 44998.5
$ws.Range("L109").Value = 44998.5
$ws.Range("N109").Value = -47772.5
$ws.Range("H113").Value = 1969.4546
$ws.Range("J113").Value = 2781.5
$ws.Range("L113").Value = 2781.5
$ws.Range("N113").Value = -7121.5
$ws.Range("H122").Value = 3238.889
$ws.Range("I122").Value = 3060.8333
$ws.Range("J122").Value = 3595
$ws.Range("K122").Value = 9182.499899999999
$ws.Range("L122").Value = 10785
$ws.Range("M122").Value = -6732.499899999999
$ws.Range("N122").Value = -15685
$ws.Range("H123").Value = 19999
$ws.Range("J123").Value = 19999
$ws.Range("L123").Value = 19999
$ws.Range("N123").Value = -29799
$ws.Range("H126").Value = 2998.3333
$ws.Range("I126").Value = 2998.5
$ws.Range("K126").Value = 8995.5
$ws.Range("M126").Value = -6525.5
$ws.Range("H132").Value = 4870.4287
$ws.Range("I132").Value = 4682.3335
$ws.Range("K132").Value = 14047.0005
$ws.Range("M132").Value = -11517.0005
$ws.Range("H136").Value = 41669584
$ws.Range("I136").Value = 3023.3
$ws.Range("K136").Value = 9069.900000000001
$ws.Range("M136").Value = -6519.900000000001

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 28675
$ws.Range("J45").Value = 29900
$ws.Range("L45").Value = 29900
$ws.Range("N45").Value = -30882
$ws.Range("H81").Value = 3949.25
$ws.Range("I81").Value = 3932.3333
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 7864.6666
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -6803.6666
$ws.Range("N81").Value = -10122
$ws.Range("H84").Value = 3949.25
$ws.Range("I84").Value = 3932.3333
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 39323.333
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -34019.333
$ws.Range("N84").Value = -50608
$ws.Range("H96").Value = 9167.166999999999
$ws.Range("I96").Value = 9334.333000000001
$ws.Range("J96").Value = 9000
$ws.Range("K96").Value = 9334.333000000001
$ws.Range("L96").Value = 9000
$ws.Range("M96").Value = -7961.333000000001
$ws.Range("N96").Value = -11746
$ws.Range("H113").Value = 1047.1428
$ws.Range("I113").Value = 982.82355
$ws.Range("K113").Value = 2948.47065
$ws.Range("M113").Value = -778.4706499999998
$ws.Range("H122").Value = 2719.2942
$ws.Range("I122").Value = 2659.2856
$ws.Range("J122").Value = 2999.3333
$ws.Range("K122").Value = 7977.8568
$ws.Range("L122").Value = 8997.999899999999
$ws.Range("M122").Value = -5527.8568
$ws.Range("N122").Value = -13897.9999
$ws.Range("H126").Value = 3194.6
$ws.Range("J126").Value = 3495
$ws.Range("L126").Value = 10485
$ws.Range("N126").Value = -15425
$ws.Range("H132").Value = 1317.8
$ws.Range("I132").Value = 1276
$ws.Range("K132").Value = 3828
$ws.Range("M132").Value = -1298
$ws.Range("H136").Value = 1695.9333
$ws.Range("I136").Value = 1555.12
$ws.Range("K136").Value = 4665.36
$ws.Range("M136").Value = -2115.36
